$d = $word.ActiveDocument

# -- Replace the -LE"..." linker-path literal with the (x86) variant --
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "-LE""c:\program files\borland\delphi7\Projects\Bpl""",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
if ($found1) {
    $rng1.Text = "-LE""c:\program files (x86)\borland\delphi7\Projects\Bpl"""
}

# -- Replace the -LN"..." linker-path literal with the (x86) variant --
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "-LN""c:\program files\borland\delphi7\Projects\Bpl""",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
if ($found2) {
    $rng2.Text = "-LN""c:\program files (x86)\borland\delphi7\Projects\Bpl"""
}
